$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: Gen -> MaxFES
$ws.Range("A1").Value = "MaxFES"

# Update column A values (rows 2-14)
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# Update AZ column values (rows 2-14) to new values before removing BA
$ws.Range("AZ2").Value = 260354986.4492386
$ws.Range("AZ3").Value = 67564349.33899023
$ws.Range("AZ4").Value = 18958509.68240635
$ws.Range("AZ5").Value = 15012133.01582103
$ws.Range("AZ6").Value = 14823030.82596275
$ws.Range("AZ7").Value = 14823030.82596058
$ws.Range("AZ8").Value = 14823030.82596058
$ws.Range("AZ9").Value = 14823030.82596058
$ws.Range("AZ10").Value = 14823030.82596058
$ws.Range("AZ11").Value = 14823030.82596058
$ws.Range("AZ12").Value = 14823030.82596058
$ws.Range("AZ13").Value = 14823030.82596058
$ws.Range("AZ14").Value = 14823030.82596058

# Delete column BA entirely (removes Run 50 header and its data)
$ws.Range("BA:BA").Delete()
